# Start of the csv-data import rework: the "Stand per ..." caption in row 3
# moves one column to the right (it now sits over F:G instead of E:F) so the
# new import column can be inserted in front of it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the caption text currently shown in the merged E3:F3 block.
$standValue = $ws.Range("E3").Value()

# Break the old E3:F3 merge before re-merging the new range.
$ws.Range("E3:F3").UnMerge()

# F3 becomes the anchor of the new merged block, carrying the same caption
# text that used to live in E3 ("Stand per 09.02.2018").
$ws.Range("F3").Value = $standValue

# Re-merge the caption over F3:G3.
$ws.Range("F3:G3").Merge()

# Land the selection back on A1, like a freshly (re)opened sheet.
$ws.Range("A1").Select()
